$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 360; this shifts the existing rows 360-394
# down to 361-395 (and brings their formatting along automatically).
$ws.Rows.Item(360).Insert()

# Populate the newly inserted row 360 with the new record's data.
$ws.Cells.Item(360, 1).Value2 = 3
$ws.Cells.Item(360, 2).Value2 = "Femacal de La Calera"
$ws.Cells.Item(360, 3).Value2 = "Coquimbo"
$ws.Cells.Item(360, 4).Value2 = 44769
$ws.Cells.Item(360, 5).Value2 = 5
$ws.Cells.Item(360, 6).Value2 = 100112040
$ws.Cells.Item(360, 7).Value2 = "Cilantro"
$ws.Cells.Item(360, 8).Value2 = "Sin especificar"
$ws.Cells.Item(360, 9).Value2 = "Primera"
$ws.Cells.Item(360, 10).Value2 = 210
$ws.Cells.Item(360, 11).Value2 = 4000
$ws.Cells.Item(360, 12).Value2 = 4500
$ws.Cells.Item(360, 13).Value2 = 4238
$ws.Cells.Item(360, 14).Value2 = "`$/docena de atados (3 kilos)"
$ws.Cells.Item(360, 15).Value2 = "Provincia de Quillota"
$ws.Cells.Item(360, 16).Value2 = 1413
$ws.Cells.Item(360, 17).Value2 = 3
$ws.Cells.Item(360, 18).Value2 = "Hortaliza"
